$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new to-do item in A21
$ws.Range("A21").Value = "read up on and have a stance to the flattening of the PC literature"

# Highlight the relevant cells (yellow fill + wrap text)
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A4").WrapText = $true

$ws.Range("A18").Interior.Color = 65535
$ws.Range("A18").WrapText = $true

$ws.Range("A20").Interior.Color = 65535
$ws.Range("A20").WrapText = $true

$ws.Range("A21").Interior.Color = 65535
$ws.Range("A21").WrapText = $true

# Update selection to mimic the cursor having moved to A29
$ws.Range("A29").Select()
